# Add retry in case of server error
# - Insert a new localization row (ServerErrorFailure) right before the
#   "ConfirmNumerousRequests" row on the Localization sheet, shifting all
#   following rows down by one and growing the Table13 table accordingly.
# - Make "Localization" the active sheet/tab instead of "Settings".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Insert a new blank row above what is currently row 43 (the
# "ConfirmNumerousRequests" entry), shifting existing rows down.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row with the new localized strings.
$ws.Range("A43").Value = "ServerErrorFailure"
$ws.Range("B43").Value = "HTTP Request failed due to server error issues."
$ws.Range("C43").Value = "サーバーエラーの問題のため、リクエストが失敗しました。"

# Match the wrap-text formatting used by the rest of column B/C in the table.
$ws.Range("B43").WrapText = $true
$ws.Range("C43").WrapText = $true

# Grow the table (Table13) so it now covers the extra row.
$lo = $ws.ListObjects.Item(1)
$lastRow = $lo.Range.Row + $lo.Range.Rows.Count - 1 + 1
$lastCol = $lo.Range.Column + $lo.Range.Columns.Count - 1
$newRange = $ws.Range($ws.Cells.Item($lo.Range.Row, $lo.Range.Column), $ws.Cells.Item($lastRow, $lastCol))
$lo.Resize($newRange)

# Switch the active tab from "Settings" to "Localization".
$ws.Activate()
